$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")
Write-Host $ws.Name
